$d = $word.ActiveDocument

# 1. Insert "This is " before the run containing "fred"
$d.Content.Find.Execute("fred", $false, $false, $false, $false, $false,
                         $true, 1, $false, "This is fred", 2)

# 2. Update overflowPunct (HangingPunctuation) on both comments from True to False
foreach ($c in $d.Comments) {
    $c.Range.ParagraphFormat.HangingPunctuation = $false
}
